$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 432, shifting existing data (rows 432-470) down to 434-472.
$ws.Rows("432:433").Insert()

# Fill in the two newly inserted rows with data.
# Row 432 (new)
$ws.Range("A432").Value = 7
$ws.Range("B432").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C432").Value = "Ñuble"
$ws.Range("D432").Value = 45106
$ws.Range("D432").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E432").Value = 16
$ws.Range("F432").Value = 100114013
$ws.Range("G432").Value = "Zanahoria"
$ws.Range("H432").Value = "Sin especificar"
$ws.Range("I432").Value = "Primera"
$ws.Range("J432").Value = 100
$ws.Range("K432").Value = 7000
$ws.Range("L432").Value = 7000
$ws.Range("M432").Value = 7000
$ws.Range("N432").Value = "$/saco 20 kilos"
$ws.Range("O432").Value = "Provincia de Diguillín"
$ws.Range("P432").Value = 350
$ws.Range("Q432").Value = 20
$ws.Range("R432").Value = "Hortaliza"

# Row 433 (new)
$ws.Range("A433").Value = 7
$ws.Range("B433").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C433").Value = "Ñuble"
$ws.Range("D433").Value = 45106
$ws.Range("D433").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E433").Value = 16
$ws.Range("F433").Value = 100114013
$ws.Range("G433").Value = "Zanahoria"
$ws.Range("H433").Value = "Sin especificar"
$ws.Range("I433").Value = "Segunda"
$ws.Range("J433").Value = 100
$ws.Range("K433").Value = 6000
$ws.Range("L433").Value = 6000
$ws.Range("M433").Value = 6000
$ws.Range("N433").Value = "$/saco 20 kilos"
$ws.Range("O433").Value = "Provincia de Diguillín"
$ws.Range("P433").Value = 300
$ws.Range("Q433").Value = 20
$ws.Range("R433").Value = "Hortaliza"
